$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting (values look numeric, e.g. "56.406.53", "1.09")
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "56.406.53"
$ws.Range("E2").Value = "  -2.33%  "

# Row 3
$ws.Range("D3").Value = "2.378.49"
$ws.Range("E3").Value = "  -3.32%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "505.50"
$ws.Range("E5").Value = "  -1.42%  "

# Row 6
$ws.Range("D6").Value = "130.21"
$ws.Range("E6").Value = "  -2.87%  "

# Row 7
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  -2.37%  "

# Row 9
$ws.Range("D9").Value = "2.388.45"
$ws.Range("E9").Value = "  -2.85%  "

# Row 10
$ws.Range("D10").Value = "0.0988"
$ws.Range("E10").Value = "  +0.89%  "

# Row 11
$ws.Range("E11").Value = "  +0.13%  "

# Row 12
$ws.Range("D12").Value = "4.84"
$ws.Range("E12").Value = "  +4.28%  "

# Row 13
$ws.Range("D13").Value = "0.327"
$ws.Range("E13").Value = "  +0.68%  "

# Row 14
$ws.Range("D14").Value = "2.801.49"
$ws.Range("E14").Value = "  -3.05%  "

# Row 15
$ws.Range("D15").Value = "56.351.40"
$ws.Range("E15").Value = "  -2.33%  "

# Row 16
$ws.Range("D16").Value = "21.51"
$ws.Range("E16").Value = "  -2.74%  "

# Row 17
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  -1.68%  "

# Row 18
$ws.Range("D18").Value = "2.371.27"
$ws.Range("E18").Value = "  -1.25%  "

# Row 19
$ws.Range("D19").Value = "10.06"
$ws.Range("E19").Value = "  -2.98%  "

# Row 20
$ws.Range("D20").Value = "4.05"
$ws.Range("E20").Value = "  -1.29%  "

# Row 21
$ws.Range("D21").Value = "309.33"
$ws.Range("E21").Value = "  -2.45%  "

# Row 22
$ws.Range("D22").Value = "6.29"
$ws.Range("E22").Value = "  -2.49%  "

# Row 23
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.22%  "

# Row 24
$ws.Range("D24").Value = "66.27"
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.29%  "

# Row 26
$ws.Range("D26").Value = "0.369"
$ws.Range("E26").Value = "  -4.08%  "

# Row 27
$ws.Range("D27").Value = "0.147"
$ws.Range("E27").Value = "  -5.76%  "

# Row 28
$ws.Range("D28").Value = "7.27"
$ws.Range("E28").Value = "  -4.77%  "

# Row 29
$ws.Range("D29").Value = "172.78"
$ws.Range("E29").Value = "  +1.18%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0713"
$ws.Range("E30").Value = "  -3.43%  "

# Row 31
$ws.Range("D31").Value = "1.65"
$ws.Range("E31").Value = "  -2.97%  "

# Row 32
$ws.Range("E32").Value = "  +0.06%  "

# Row 33
$ws.Range("D33").Value = "5.83"
$ws.Range("E33").Value = "  -4.39%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.09"
$ws.Range("E34").Value = "  -5.15%  "

# Row 35
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$ws.Range("D36").Value = "17.72"
$ws.Range("E36").Value = "  -2.21%  "

# Row 37
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  -2.95%  "

# Row 38
$ws.Range("D38").Value = "3.75"
$ws.Range("E38").Value = "  -3.89%  "

# Row 39
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").Value = "0.820"
$ws.Range("E39").Value = "  +0.80%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "36.51"
$ws.Range("E40").Value = "  -0.85%  "

# Row 41
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  -5.52%  "

# Row 42
$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  -0.85%  "

# Row 43
$ws.Range("D43").Value = "128.88"
$ws.Range("E43").Value = "  -3.83%  "

# Row 44
$ws.Range("D44").Value = "4.93"
$ws.Range("E44").Value = "  -2.88%  "

# Row 45
$ws.Range("D45").Value = "0.570"
$ws.Range("E45").Value = "  -0.68%  "

# Row 46
$ws.Range("D46").Value = "0.0897"
$ws.Range("E46").Value = "  -2.15%  "

# Row 47
$ws.Range("D47").Value = "240.52"
$ws.Range("E47").Value = "  -6.53%  "

# Row 48
$ws.Range("D48").Value = "0.0482"
$ws.Range("E48").Value = "  -2.54%  "

# Row 49
$ws.Range("D49").Value = "0.0207"
$ws.Range("E49").Value = "  -2.96%  "

# Row 50
$ws.Range("D50").Value = "17.17"
$ws.Range("E50").Value = "  -1.34%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "1.54"
$ws.Range("E51").Value = "  -4.55%  "
